# Queues implemented... some invoking done, will finish after i deploy to
# orchestrator... invoice builds correctly... non-fatal error writing to
# Orders2.xlsx
#
# Update the "Address" sheet's billing/shipping example data, then leave the
# "Address" sheet as the active tab/selection (matching the author's last
# on-screen state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Address")

# --- Shipping address block (row 5) -> "William Gates", Microsoft Street ---
$ws.Range("B5").Value = "William"
$ws.Range("C5").Value = "Gates"
$ws.Range("D5").Value = "Microsoft Street"
$ws.Range("F5").Value = "MicroSoft"

# --- Billing address block (row 2) -> new street + phone numbers ---
$ws.Range("D2").Value = "Apple Avenue"
$ws.Range("F2").Value = 789
$ws.Range("G2").Value = 55511234
$ws.Range("H2").Value = "Apple"

# Column D widened slightly to fit the new street text.
$ws.Columns.Item(4).ColumnWidth = 14.167389418907199

# The author ended their session with the "Address" sheet active and cell
# G2 selected (previously "Orders" was the active/selected tab).
$ws.Activate() | Out-Null
$ws.Range("G2").Select() | Out-Null
